$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Send_status" column (C) header, styled like the existing A1/B1 headers
$ws.Range("C1").Value = "Send_status"
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)

# New boolean (FALSE) values for the data rows in column C
$ws.Range("C2").Value = $false
$ws.Range("C3").Value = $false
$ws.Range("C4").Value = $false
